{"js": "// Budget justification RL.docx \u2014 add data-infrastructure approaches text.\n//\n// The only substantive content change in the target revision (the rest of\n// the diff is Word's own run-splitting / proofing-mark churn) is inside the\n// \"Firstly, cloud infrastructure...\" paragraph, where the list of specific\n// cloud services was expanded:\n//\n//   ...include: Azure Blob storage ($0 - $1000), Azure static web app...\n//\n// becomes\n//\n//   ...include: Azure Blob/Data-Lake storage ($0 - $1000), Azure Synapse\n//   storage ($0 - $1000),  Azure static web app...\n//\n// i.e. two insertions:\n//   1) \"/Data-Lake\" right after \"Azure Blob\"\n//   2) \" Azure Synapse storage ($0 - $1000), \" right after\n//      \"storage ($0 - $1000),\"\n\nconst body = context.document.body;\n\n// 1) \"Azure Blob\" -> \"Azure Blob/Data-Lake\"\nlet blobResults = body.search(\"Azure Blob\", { matchCase: true });\nblobResults.load(\"items\");\nawait context.sync();\n\nif (blobResults.items.length > 0) {\n  blobResults.items[0].insertText(\"/Data-Lake\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// 2) \"storage ($0 - $1000),\" -> \"storage ($0 - $1000), Azure Synapse storage ($0 - $1000), \"\nlet storageResults = body.search(\"storage ($0 - $1000),\", { matchCase: true });\nstorageResults.load(\"items\");\nawait context.sync();\n\nif (storageResults.items.length > 0) {\n  storageResults.items[0].insertText(\n    \" Azure Synapse storage ($0 - $1000), \",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n}\n", "ps1": "# Budget justification RL.docx \u2014 add data-infrastructure approaches text.\n#\n# The only substantive content change in the target revision (the rest of\n# the diff is Word's own run-splitting / proofing-mark churn) is inside the\n# \"Firstly, cloud infrastructure...\" paragraph, where the list of specific\n# cloud services was expanded:\n#\n#   ...include: Azure Blob storage ($0 - $1000), Azure static web app...\n#\n# becomes\n#\n#   ...include: Azure Blob/Data-Lake storage ($0 - $1000), Azure Synapse\n#   storage ($0 - $1000),  Azure static web app...\n#\n# i.e. two insertions:\n#   1) \"/Data-Lake\" right after \"Azure Blob\"\n#   2) \" Azure Synapse storage ($0 - $1000), \" right after\n#      \"storage ($0 - $1000),\"\n\n$d = $word.ActiveDocument\n\n# 1) \"Azure Blob\" -> \"Azure Blob/Data-Lake\"\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWildcards = $false\n$rng.Find.Text = \"Azure Blob\"\nif ($rng.Find.Execute()) {\n    $rng.Collapse($wdCollapseEnd)\n    $rng.InsertAfter(\"/Data-Lake\")\n}\n\n# 2) \"storage ($0 - $1000),\" -> \"storage ($0 - $1000), Azure Synapse storage ($0 - $1000), \"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.MatchCase = $true\n$rng2.Find.MatchWildcards = $false\n$rng2.Find.Text = \"storage (`$0 - `$1000),\"\nif ($rng2.Find.Execute()) {\n    $rng2.Collapse($wdCollapseEnd)\n    $rng2.InsertAfter(\" Azure Synapse storage (`$0 - `$1000), \")\n}\n"}
